$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.876.33"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "2.270.20"
$ws.Range("E3").Value = "  -3.62%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.571"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.62%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.505"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0798"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.11%  "
$ws.Range("E13").Value = "  -1.74%  "
$ws.Range("D14").Value = "2.614.48"
$ws.Range("E14").Value = "  -3.66%  "
$ws.Range("D15").Value = "2.266.39"
$ws.Range("E15").Value = "  -3.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.16%  "
$ws.Range("D17").Value = "46.858.20"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.796"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.60%  "
$ws.Range("E19").Value = "  +1.99%  "
$ws.Range("E20").Value = "  -8.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "246.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("E24").Value = "  -7.72%  "
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.76%  "
$ws.Range("E27").Value = "  -2.36%  "
$ws.Range("E28").Value = "  -3.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("E31").Value = "  +7.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "145.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0766"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.32%  "
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("E37").Value = "  -2.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "15.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.26%  "
$ws.Range("E39").Value = "  -10.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0296"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -11.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "94.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +15.76%  "
$ws.Range("D45").Value = "1.780.72"
$ws.Range("E45").Value = "  -5.13%  "
$ws.Range("E46").Value = "  -4.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "70.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.41%  "
$ws.Range("E48").Value = "  -8.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "94.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.19%  "
